$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the value for C2 and C3 to match the existing shared string
# "C:\Program Files (x86)" already used in C4
$ws.Range("C2").Value = "C:\Program Files (x86)"
$ws.Range("C3").Value = "C:\Program Files (x86)"

# Update the active selection to C2 (matches the diff's <selection activeCell="C2" sqref="C2"/>)
$ws.Range("C2").Select()
